# Generate Report for Handback
# Updates the localization-status workbook after a successful handback:
#  - Overview sheet: status text changes from "Ready for handoff" to
#    "Handed back: in sync with en-US" and the two status columns widen.
#  - zh-cn / de-de sheets: the "Latest Handback DateTime" is refreshed,
#    the stale "Error Detail" message is cleared (handback is now in
#    sync with en-US so there's no error left to report), and the
#    Status / Error Detail columns are resized to fit the new content.

$wb = $excel.ActiveWorkbook

# ---- Overview sheet ----
$overview = $wb.Worksheets.Item("Overview")
$overview.Range("E2").Value = "Handed back: in sync with en-US"
$overview.Range("F2").Value = "Handed back: in sync with en-US"
$overview.Columns.Item(5).ColumnWidth = 29.9777050018311
$overview.Columns.Item(6).ColumnWidth = 29.9777050018311

# ---- zh-cn sheet ----
$zhcn = $wb.Worksheets.Item("zh-cn")
$zhcn.Range("K2").Value = "2016-08-30 08:47:01"
$zhcn.Range("P2").Value = ""
$zhcn.Columns.Item(3).ColumnWidth = 29.9777050018311
$zhcn.Columns.Item(16).ColumnWidth = 13.7470531463623

# ---- de-de sheet ----
$dede = $wb.Worksheets.Item("de-de")
$dede.Range("K2").Value = "2016-08-30 08:47:27"
$dede.Range("P2").Value = ""
$dede.Columns.Item(3).ColumnWidth = 29.9777050018311
$dede.Columns.Item(16).ColumnWidth = 13.7470531463623
